$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "22.443.98"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.573.26"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.04%  "
Set-TextValue $ws.Range("D6") "291.35"
$ws.Range("E6").Value = "  +0.02%  "
Set-TextValue $ws.Range("D7") "0.3741"
$ws.Range("E7").Value = "  -0.80%  "
Set-TextValue $ws.Range("D8") "49.87"
$ws.Range("E8").Value = "  +0.11%  "
Set-TextValue $ws.Range("D9") "0.3399"
$ws.Range("E9").Value = "  -0.84%  "
Set-TextValue $ws.Range("D10") "0.07555"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("E12").Value = "  +0.03%  "
Set-TextValue $ws.Range("D13") "21.35"
$ws.Range("E13").Value = "  +0.54%  "
Set-TextValue $ws.Range("D14") "5.993"
$ws.Range("E14").Value = "  -0.27%  "
Set-TextValue $ws.Range("D15") "6.931"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "1.566.74"
$ws.Range("E16").Value = "  -0.26%  "
Set-TextValue $ws.Range("D18") "91.02"
$ws.Range("E18").Value = "  +0.80%  "
Set-TextValue $ws.Range("D19") "0.06737"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("E20").Value = "  +0.03%  "
Set-TextValue $ws.Range("D21") "6.259"
$ws.Range("E21").Value = "  +0.57%  "
Set-TextValue $ws.Range("D22") "16.40"
$ws.Range("E22").Value = "  -2.49%  "
Set-TextValue $ws.Range("D23") "12.14"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "22.447.33"
$ws.Range("E24").Value = "  +0.07%  "
Set-TextValue $ws.Range("D25") "2.332"
$ws.Range("E25").Value = "  -3.99%  "
Set-TextValue $ws.Range("D26") "2.607"
$ws.Range("E26").Value = "  -4.60%  "
Set-TextValue $ws.Range("D27") "20.12"
$ws.Range("E27").Value = "  -0.89%  "
Set-TextValue $ws.Range("D29") "5.002"
$ws.Range("E29").Value = "  -0.61%  "
Set-TextValue $ws.Range("D30") "125.61"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "1.747.16"
$ws.Range("E31").Value = "  +0.09%  "
Set-TextValue $ws.Range("D32") "1.052"
$ws.Range("E32").Value = "  +5.09%  "
Set-TextValue $ws.Range("D33") "6.119"
$ws.Range("E33").Value = "  -1.41%  "
Set-TextValue $ws.Range("D34") "1.982"
Set-TextValue $ws.Range("D35") "9.821"
$ws.Range("E35").Value = "  -2.12%  "
Set-TextValue $ws.Range("D36") "0.08412"
$ws.Range("E36").Value = "  -1.85%  "
$ws.Range("E37").Value = "  +3.31%  "
Set-TextValue $ws.Range("D38") "0.02463"
$ws.Range("E38").Value = "  -3.32%  "
Set-TextValue $ws.Range("D39") "0.2286"
$ws.Range("E39").Value = "  -1.42%  "
Set-TextValue $ws.Range("D40") "0.06511"
$ws.Range("E40").Value = "  -1.14%  "
Set-TextValue $ws.Range("D41") "5.453"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("E42").Value = "  -2.66%  "
Set-TextValue $ws.Range("D43") "0.6238"
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("E44").Value = "  +0.04%  "
Set-TextValue $ws.Range("D45") "13.92"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("E46").Value = "  +0.30%  "
Set-TextValue $ws.Range("D47") "0.5815"
$ws.Range("E47").Value = "  -3.37%  "
Set-TextValue $ws.Range("D50") "1.221"
$ws.Range("E50").Value = "  -6.56%  "
Set-TextValue $ws.Range("D51") "0.07320"
$ws.Range("E51").Value = "  -0.17%  "

# Row 48/49 swap: Quant <-> NEARProtocol
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D48") "2.081"
$ws.Range("E48").Value = "  -0.30%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D49") "129.32"
$ws.Range("E49").Value = "  +3.22%  "
